$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.268.01'
$ws.Range("E2").Value = '  +5.80%  '
$ws.Range("D3").Value = '2.788.89'
$ws.Range("E3").Value = '  +6.22%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'116.60"
$ws.Range("E5").Value = '  +4.37%  '
$ws.Range("D6").Value = "'340.50"
$ws.Range("E6").Value = '  +4.59%  '
$ws.Range("D7").Value = "'0.554"
$ws.Range("E7").Value = '  +5.60%  '
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +5.84%  '
$ws.Range("D10").Value = "'42.02"
$ws.Range("E10").Value = '  +6.61%  '
$ws.Range("E11").Value = '  +6.94%  '
$ws.Range("D12").Value = "'20.08"
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("D14").Value = "'7.63"
$ws.Range("E14").Value = '  +1.12%  '
$ws.Range("D15").Value = '3.231.33'
$ws.Range("E15").Value = '  +6.45%  '
$ws.Range("D16").Value = '2.799.60'
$ws.Range("E16").Value = '  +6.64%  '
$ws.Range("D17").Value = "'0.885"
$ws.Range("E17").Value = '  +3.93%  '
$ws.Range("D18").Value = '52.070.92'
$ws.Range("E18").Value = '  +5.51%  '
$ws.Range("D19").Value = "'3.20"
$ws.Range("E19").Value = '  +10.22%  '
$ws.Range("D20").Value = "'13.34"
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D22").Value = '0.0₃0982'
$ws.Range("E22").Value = '  +3.75%  '
$ws.Range("D23").Value = "'278.43"
$ws.Range("E23").Value = '  +3.86%  '
$ws.Range("D24").Value = "'70.22"
$ws.Range("E24").Value = '  +1.85%  '
$ws.Range("D25").Value = "'2.78"
$ws.Range("E25").Value = '  +8.83%  '
$ws.Range("D26").Value = "'26.86"
$ws.Range("E26").Value = '  +3.46%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").Value = "'10.27"
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = '  +1.30%  '
$ws.Range("D30").Value = "'0.143"
$ws.Range("E30").Value = '  +3.88%  '
$ws.Range("D31").Value = "'34.82"
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("D32").Value = "'50.45"
$ws.Range("E32").Value = '  +1.88%  '
$ws.Range("D33").Value = "'5.75"
$ws.Range("E33").Value = '  +5.25%  '
$ws.Range("D34").Value = "'0.0829"
$ws.Range("E34").Value = '  +2.70%  '
$ws.Range("E35").Value = '  +5.02%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").Value = "'4.95"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("E39").Value = '  +4.95%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = "'0.0378"
$ws.Range("E40").Value = '  +13.80%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = "'2.72"
$ws.Range("E41").Value = '  +27.74%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = "'0.116"
$ws.Range("E42").Value = '  +4.34%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = "'23.32"
$ws.Range("E43").Value = '  +3.45%  '
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("D45").Value = "'124.89"
$ws.Range("E45").Value = '  -3.74%  '
$ws.Range("D46").Value = '2.098.28'
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("E47").Value = '  +1.90%  '
$ws.Range("E48").Value = '  +3.51%  '
$ws.Range("D49").Value = "'5.58"
$ws.Range("E49").Value = '  +7.48%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = "'8.99"
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("B51").Value = 'SEI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").Value = "'0.893"
$ws.Range("E51").Value = '  +20.36%  '
